# Weekly update: a new week of Brócoli price data (fecha 44466) was added
# for "Terminal La Palmera de La Serena", inserted right before the
# existing 44389 week (rows 269/270), pushing every subsequent row down
# by two. The two freshly inserted rows reuse the same Primera/Segunda
# price figures that used to sit in (old) rows 269/270 — only the date
# changes for them; everything else cascades down unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right above the current row 269 (Primera /
# Segunda pair for fecha 44389). Excel shifts rows 269.. down to 271..,
# growing the used range to A1:R390 automatically.
$ws.Rows("269:270").Insert()

# The data that used to be in rows 269/270 is now in rows 271/272.
# Duplicate it back up into the freshly inserted 269/270 rows so both
# weeks carry identical Primera/Segunda pricing...
$ws.Range("A271:R272").Copy($ws.Range("A269:R270"))

# ...except the new rows represent a newer market date.
$ws.Range("D269:D270").Value = 44466
